# Add "2022-Q3" data: insert a new quarter worksheet + a new summary row.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    sheet that currently holds the "2022-Q2" data), so every following
#    quarter sheet shifts down by one tab position - matching the diff's
#    sheetId / r:id renumbering.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$firstQuarterSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($firstQuarterSheet)
$q3.Name = "2022-Q3"

# Header style (bordered / bold, centered) matches the other quarter sheets'
# row-1 + column-A cells, which all use the same cell style. Grab it from the
# sheet we just pushed down instead of hard-coding a style index.
$firstQuarterSheet.Cells.Item(1, 2).Copy() | Out-Null
$q3.Cells.Item(1, 2).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(1, 3).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(1, 5).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(1, 7).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(1, 8).PasteSpecial(-4122) | Out-Null

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

$firstQuarterSheet.Cells.Item(2, 1).Copy() | Out-Null
$q3.Cells.Item(2, 1).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$q3.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null

# Force text storage (matches the source data, which keeps these numeric-
# looking figures and fund codes as text) for columns B-G, rows 2-4.
$q3.Range("B2:G4").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "002423"
$q3.Cells.Item(2, 3).Value = "华宝标普美国品质消费股票（LOF）美元"
$q3.Cells.Item(2, 4).Value = "3.59"
$q3.Cells.Item(2, 5).Value = "94.37"
$q3.Cells.Item(2, 6).Value = "4.34"
$q3.Cells.Item(2, 7).Value = "0.1558"
$q3.Cells.Item(2, 8).Value = 4

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "162415"
$q3.Cells.Item(3, 3).Value = "华宝标普美国品质消费股票（LOF）人民币A"
$q3.Cells.Item(3, 4).Value = "2.86"
$q3.Cells.Item(3, 5).Value = "94.37"
$q3.Cells.Item(3, 6).Value = "4.34"
$q3.Cells.Item(3, 7).Value = "0.1241"
$q3.Cells.Item(3, 8).Value = 4

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "009975"
$q3.Cells.Item(4, 3).Value = "华宝标普美国品质消费股票（LOF）人民币C"
$q3.Cells.Item(4, 4).Value = "0.73"
$q3.Cells.Item(4, 5).Value = "94.37"
$q3.Cells.Item(4, 6).Value = "4.34"
$q3.Cells.Item(4, 7).Value = "0.0317"
$q3.Cells.Item(4, 8).Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the 2022-Q3 row at row 2 and
#    shift the previously-existing rows down by one, re-numbering the
#    leading index column (A) sequentially.
# ---------------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q3", 3, 0.31),
    @("2022-Q2", 5, 0.31),
    @("2022-Q1", 3, 0.33),
    @("2021-Q4", 3, 0.38),
    @("2021-Q3", 3, 0.32),
    @("2021-Q2", 3, 0.28),
    @("2021-Q1", 3, 0.26),
    @("2020-Q4", 3, 0.22)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $row = $i + 2
    $vals = $summaryRows[$i]
    $totalSheet.Cells.Item($row, 1).Value = $i
    $totalSheet.Cells.Item($row, 2).Value = $vals[0]
    $totalSheet.Cells.Item($row, 3).Value = $vals[1]
    $totalSheet.Cells.Item($row, 4).Value = $vals[2]
}

# Style the newly-added row 9 index cell (A9) like the others in column A.
$totalSheet.Cells.Item(8, 1).Copy() | Out-Null
$totalSheet.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$totalSheet.Cells.Item(9, 1).Value = 7

# Keep "总计" as the active tab (matches the unchanged bookViews/activeTab).
$totalSheet.Activate()
